$d = $word.ActiveDocument

# The "Test Function" paragraph currently spells out the signature across
# several differently-formatted runs:
#   "int" + " validatePackageWeight(" + "double" + " " + "num" + ")"
# The edit collapses all of that into a single run (keeping the formatting
# of the original "int" run) containing the text:
#   "int validatePackageWeight(double);"

$find = $d.Content
$find.Find.ClearFormatting()
$found = $find.Find.Execute("int validatePackageWeight(double num)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $find.Start
    $matchEnd = $find.End

    # Remove everything after the word "int" (this is the part coming from
    # the other 5 runs: " validatePackageWeight(", "double", " ", "num", ")").
    $tail = $d.Range($matchStart + 3, $matchEnd)
    $tail.Text = ""

    # Replace the remaining "int" run's own text in place so the run keeps
    # its original formatting (Cascadia Mono, blue, size 19, en-HK) instead
    # of spawning a freshly-formatted run.
    $introRun = $d.Range($matchStart, $matchStart + 3)
    $introRun.Text = "int validatePackageWeight(double);"
}
